$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Item 7 ("Entrar a mostrar estudiante y despues entrar a mostrar profesor")
# is no longer applicable -> clear its description text.
$ws.Range("B9").Value = ""

# Item 1 ("Hacer que en las tablas...") is now checked off -> add the
# checkmark indicator in column P, matching the look of the other
# checked item (copy its format, then bump size/weight).
$ws.Range("P6").Copy()
$ws.Range("P3").PasteSpecial(-4122)
$ws.Range("P3").Value = "      ✔"
$ws.Range("P3").Font.Size = 12
$ws.Range("P3").Font.Bold = $true

$ws.Range("L14").Select()
